$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add two new salespeople under the same store (AL357) as row 12 (LUIS CHÁVEZ)
$ws.Range("A13").Value = "AL357"
$ws.Range("B13").Value = "JOSE CRUZ"

$ws.Range("A14").Value = "AL357"
$ws.Range("B14").Value = "JOSE MARTÍNEZ"

# Update the selection to match the saved state
$ws.Range("B22").Select()
